{"js": "// The document body underwent an overall \"accept changes\" style cleanup\n// (spelling/grammar proofing marks removed and adjacent runs merged) that\n// does not alter the visible text anywhere except for one spot: the list\n// of topics read about was expanded from bare terms to fuller phrases.\n//\n// \"gaussian,sobel,scharr,otsu\"\n//   -> \"gaussian blur,sobel edge detection,scharr edge detection,otsu thresholding\"\n\nconst searchText = \"gaussian,sobel,scharr,otsu\";\nconst replacementText =\n  \"gaussian blur,sobel edge detection,scharr edge detection,otsu thresholding\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document underwent an overall \"accept changes\" style cleanup\n# (spelling/grammar proofing marks removed and adjacent runs merged) that\n# does not alter the visible text anywhere except for one spot: the list\n# of topics read about was expanded from bare terms to fuller phrases.\n#\n# \"gaussian,sobel,scharr,otsu\"\n#   -> \"gaussian blur,sobel edge detection,scharr edge detection,otsu thresholding\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Execute(\n    \"gaussian,sobel,scharr,otsu\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"gaussian blur,sobel edge detection,scharr edge detection,otsu thresholding\",\n    2\n)\n"}
